$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the three new header cells
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting (bold, border, centered/top aligned) from an
# existing header cell onto the new header cells
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the season record (Wins/Losses/Ties) for every data row
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 85
    $ws.Cells.Item($r, 31).Value = 77
    $ws.Cells.Item($r, 32).Value = 0
}
